$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "readme" sheet: the JobNo/Date/sheet_name/Author columns were built in
#    the wrong order by the report generator; this run's fix realigns the
#    table header labels with the data that actually belongs under them and
#    refreshes the run's metadata (author, date, job folder).
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("readme")

# Corrected header order: index, Author, sheet_name, Date, JobNo
$ws.Cells.Item(1, 2).Value = "Author"
$ws.Cells.Item(1, 3).Value = "sheet_name"
$ws.Cells.Item(1, 4).Value = "Date"
$ws.Cells.Item(1, 5).Value = "JobNo"

# The "Date" column holds a numeric-looking value ("20220615"); format it as
# text first so Excel doesn't silently coerce it into a number.
$ws.Range("D2:D12").NumberFormat = "@"

$lastRow = 12
for ($r = 2; $r -le $lastRow; $r++) {
    # sheet_name now lives in column C; preserve the per-row value that used
    # to sit in column D before the header/column shuffle.
    $sheetName = $ws.Cells.Item($r, 4).Value()

    $ws.Cells.Item($r, 2).Value = "jovyan"
    $ws.Cells.Item($r, 3).Value = $sheetName
    $ws.Cells.Item($r, 4).Value = "20220615"
    $ws.Cells.Item($r, 5).Value = "/c/e"
}

# ---------------------------------------------------------------------------
# 2) "Project Information" sheet: bump the recorded analysis timestamp to
#    match this (re-)run.
# ---------------------------------------------------------------------------
$wsProj = $wb.Worksheets.Item("Project Information")
for ($r = 2; $r -le $wsProj.UsedRange.Rows.Count; $r++) {
    if ($wsProj.Cells.Item($r, 1).Value() -eq "Date of Analysis") {
        $wsProj.Cells.Item($r, 2).Value = "2022-06-15 10:33:16.688962"
    }
}
